$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: merge "Testing" " " "custom" " " "properties" runs into one run ---
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "__temp_title__"
$title.Text = "Testing custom properties"

# --- Subtitle shape: merge word-split runs into single runs, keeping the two <a:br/> ---
$subtitle = $s.Shapes.Item(2).TextFrame.TextRange

# First segment: "This" " " "is" " " "a" " " "subtitle" (18 chars) -> single run
$seg1 = $subtitle.Characters(1, 18)
$seg1.Text = "__temp_seg1__"
$subtitle.Characters(1, 13).Text = "This is a subtitle"

# Second segment (after the two <a:br/>): "A." " " "M." (5 chars) -> single run
$seg2 = $subtitle.Characters(21, 5)
$seg2.Text = "__temp_seg2__"
$subtitle.Characters(21, 13).Text = "A. M."
